$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = -215.3074555466683
$ws.Range("G3").Value = -200.1541353621739
$ws.Range("D4").Value = 85.38998748436885
$ws.Range("G4").Value = -190.371676368487
$ws.Range("D5").Value = 83.65655902870688
$ws.Range("G5").Value = -186.2626070719929
$ws.Range("G6").Value = -185.5728286247045
$ws.Range("G7").Value = -190.3167594965325
$ws.Range("G8").Value = -210.1803155584764
$ws.Range("D9").Value = 107.9386898447804
$ws.Range("F9").Value = 107.9386898447803
$ws.Range("G9").Value = -245.46707268136
$ws.Range("G10").Value = -269.3873981896476
$ws.Range("B11").Value = -432.9753410696505
$ws.Range("C11").Value = 324.7315058022379
$ws.Range("D11").Value = 121.2330954995022
$ws.Range("E11").Value = 121.2330954995021
$ws.Range("G11").Value = -279.9073060755944
$ws.Range("G12").Value = -283.8477475878325
$ws.Range("B13").Value = -437.1044095295265
$ws.Range("C13").Value = 327.8283071471449
$ws.Range("F13").Value = 122.3892346682675
$ws.Range("G13").Value = -283.0008498503894
$ws.Range("G14").Value = -280.232292859933
$ws.Range("E15").Value = 120.7165677442803
$ws.Range("F15").Value = 120.7165677442802
$ws.Range("G15").Value = -278.5311582840544
$ws.Range("B16").Value = -417.8068230501005
$ws.Range("C16").Value = 313.3551172875754
$ws.Range("D16").Value = 116.9859104540282
$ws.Range("G16").Value = -268.6933036005295
$ws.Range("B17").Value = -409.4007760748646
$ws.Range("E17").Value = 114.6322173009621
$ws.Range("G17").Value = -262.5708380849388
$ws.Range("C18").Value = 303.3601962836231
$ws.Range("G18").Value = -259.0146630298199
$ws.Range("G19").Value = -257.8044238579542
$ws.Range("G20").Value = -263.2261217359288
$ws.Range("B21").Value = -434.4988625867515
$ws.Range("C21").Value = 325.8741469400637
$ws.Range("E21").Value = 121.6596815242905
$ws.Range("G21").Value = -281.0465719734139
$ws.Range("G22").Value = -292.4481453277631
$ws.Range("B23").Value = -441.5876354132125
$ws.Range("G23").Value = -286.3814958821378
$ws.Range("G24").Value = -262.9299816402682
$ws.Range("G25").Value = -236.2929982387788
